# Dashboard testcases added and Screen shots added in HTML reports
#
# Adds three new Suite rows (DEFAULTDASHBOARD / ADDWIDGETS / REMOVEWIDGETS,
# each with Runmode "N"), moves the active selection to D4, and restores the
# saved window height.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New testcases appended below the existing CDPDASHBOARD row.
$ws.Range("A5").Value = "DEFAULTDASHBOARD"
$ws.Range("B5").Value = "N"

$ws.Range("A6").Value = "ADDWIDGETS"
$ws.Range("B6").Value = "N"

$ws.Range("A7").Value = "REMOVEWIDGETS"
$ws.Range("B7").Value = "N"

# Match the existing Runmode column formatting (center aligned, style index 5).
$ws.Range("B5:B7").HorizontalAlignment = -4108

# Selection moved from C3 to D4.
$ws.Range("D4").Select()

# Saved window geometry.
$excel.ActiveWindow.Height = 4560
